# Update "Förändrad" (changed) date column C for all existing data rows
# (rows 2-129) from 45172 -> 45175.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C129").Value = 45175

# Row 129 gains an explicit row height (15pt, custom) in the edited file.
$ws.Rows.Item(129).RowHeight = 15

# Two brand-new logging-notification rows are appended at the bottom of the
# sheet: 130 (A 41214-2023) and 131 (A 41265-2023). Columns H:Q are the
# species/threat counters (all zero for these new entries) and R is the
# (empty) wrapped "Artnamn" cell, matching every other row in the sheet.

$newRows = @(
    @{ Row = 130; Beteckning = "A 41214-2023"; Datum = 45174; Area = 51.6 },
    @{ Row = 131; Beteckning = "A 41265-2023"; Datum = 45174; Area = 29.7 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.Beteckning

    $ws.Cells.Item($row, 2).Value = $r.Datum
    $ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($row, 3).Value = 45175
    $ws.Cells.Item($row, 3).NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($row, 4).Value = "VÄRMLANDS LÄN"
    $ws.Cells.Item($row, 5).Value = "STORFORS"
    $ws.Cells.Item($row, 6).Value = "Bergvik skog väst AB"
    $ws.Cells.Item($row, 7).Value = $r.Area

    for ($col = 8; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = 0
    }

    $ws.Cells.Item($row, 18).Value = ""
    $ws.Cells.Item($row, 18).WrapText = $true
}

# Only the first new row (130) carries the explicit custom row height;
# row 131 is left at the default.
$ws.Rows.Item(130).RowHeight = 15
